# Updated borders of the board
# Flip the maze "wall" (W) / "corridor" (C) markers on a handful of cells
# to redraw the board borders, and update the view's selection/scroll
# position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O7").Value = "C"
$ws.Range("P7").Value = "C"
$ws.Range("U7").Value = "W"

$ws.Range("O8").Value = "C"
$ws.Range("P8").Value = "W"
$ws.Range("Q8").Value = "W"
$ws.Range("T8").Value = "C"
$ws.Range("U8").Value = "W"

$ws.Range("O9").Value = "W"
$ws.Range("P9").Value = "W"
$ws.Range("T9").Value = "C"
$ws.Range("U9").Value = "W"

$ws.Range("S10").Value = "C"

$ws.Range("R11").Value = "C"
$ws.Range("T11").Value = "W"

$ws.Range("S12").Value = "W"

$ws.Range("S14").Value = "W"

$ws.Range("R15").Value = "C"
$ws.Range("T15").Value = "W"

$ws.Range("S16").Value = "C"

$ws.Range("O17").Value = "W"
$ws.Range("P17").Value = "W"
$ws.Range("T17").Value = "C"
$ws.Range("U17").Value = "W"

$ws.Range("O18").Value = "C"
$ws.Range("P18").Value = "W"
$ws.Range("Q18").Value = "W"
$ws.Range("T18").Value = "C"
$ws.Range("U18").Value = "W"

$ws.Range("O19").Value = "C"
$ws.Range("P19").Value = "C"
$ws.Range("U19").Value = "W"
$ws.Range("BB19").Value = "C"

$ws.Range("T20").Value = "C"
$ws.Range("U20").Value = "W"
$ws.Range("AK20").Value = "W"
$ws.Range("AL20").Value = "C"
$ws.Range("BB20").Value = "C"

$ws.Range("H21").Value = "W"
$ws.Range("I21").Value = "W"
$ws.Range("T21").Value = "C"
$ws.Range("U21").Value = "W"
$ws.Range("AK21").Value = "W"
$ws.Range("AL21").Value = "C"
$ws.Range("BA21").Value = "W"
$ws.Range("BC21").Value = "W"

$ws.Range("I22").Value = "C"
$ws.Range("U22").Value = "W"
$ws.Range("AK22").Value = "W"
$ws.Range("BB22").Value = "C"

# Update the view: scroll/selection moved from B4/I28 to A4/D6:BG23
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("D6:BG23").Select()
